$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "ether" row/column (row 3 and column C), shrinking the
# correlation matrix from a 5x5 (A1:F6) to a 4x4 (A1:E5) table of
# functional groups.
$ws.Range("C1:C6").Delete()
$ws.Range("A3:F3").Delete()

# Rewrite the remaining header row and data values explicitly so the
# final values match exactly regardless of how the Delete() shifted
# things.
$ws.Range("B1").Value = "amines"
$ws.Range("C1").Value = "aldehyde"
$ws.Range("D1").Value = "aromatic"
$ws.Range("E1").Value = "cycle"

$ws.Range("A2").Value = "amines"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 0.08663954262881006
$ws.Range("D2").Value = -0.07991058110187403
$ws.Range("E2").Value = 0.09542827420334238

$ws.Range("A3").Value = "aldehyde"
$ws.Range("B3").Value = 0.08663954262881006
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = -0.05656442953322588
$ws.Range("E3").Value = 0.128529664443122

$ws.Range("A4").Value = "aromatic"
$ws.Range("B4").Value = -0.07991058110187403
$ws.Range("C4").Value = -0.05656442953322588
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = -0.2608495356690502

$ws.Range("A5").Value = "cycle"
$ws.Range("B5").Value = 0.09542827420334238
$ws.Range("C5").Value = 0.128529664443122
$ws.Range("D5").Value = -0.2608495356690502
$ws.Range("E5").Value = 1
